# Scheduled data-refresh: update cached market/profit figures (columns
# H:N) across several sheets. No formulas are involved anywhere in this
# workbook - every cell below is a literal, externally-computed value
# being overwritten in place, matching the upstream source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 990.1429000000001
$ws.Range("J17").Value = 991.425
$ws.Range("L17").Value = 2974.275
$ws.Range("N17").Value = -3310.275
$ws.Range("H28").Value = 562.36365
$ws.Range("I28").Value = 568.7
$ws.Range("K28").Value = 568.7
$ws.Range("M28").Value = -83.70000000000005
$ws.Range("H69").Value = 9255.117
$ws.Range("I69").Value = 6448.5
$ws.Range("K69").Value = 19345.5
$ws.Range("M69").Value = -18471.5
$ws.Range("H72").Value = 9255.117
$ws.Range("I72").Value = 6448.5
$ws.Range("K72").Value = 58036.5
$ws.Range("M72").Value = -53668.5
$ws.Range("H100").Value = 7014.125
$ws.Range("I100").Value = 2601
$ws.Range("K100").Value = 2601
$ws.Range("M100").Value = -2060
$ws.Range("H137").Value = 1048261.5
$ws.Range("I137").Value = 1058.0555
$ws.Range("K137").Value = 3174.1665
$ws.Range("M137").Value = -624.1664999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 793354.3
$ws.Range("I132").Value = 1079470.5
$ws.Range("J132").Value = 6534.875
$ws.Range("K132").Value = 3238411.5
$ws.Range("L132").Value = 19604.625
$ws.Range("M132").Value = -3235881.5
$ws.Range("N132").Value = -24664.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 785.91174
$ws.Range("I94").Value = 701.6818
$ws.Range("K94").Value = 701.6818
$ws.Range("M94").Value = -250.6818
$ws.Range("H134").Value = 657061.4
$ws.Range("I134").Value = 896517.25
$ws.Range("K134").Value = 2689551.75
$ws.Range("M134").Value = -2687016.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 83356.21000000001
$ws.Range("I31").Value = 131352.5
$ws.Range("K31").Value = 131352.5
$ws.Range("M31").Value = -131057.5
$ws.Range("H34").Value = 83356.21000000001
$ws.Range("I34").Value = 131352.5
$ws.Range("K34").Value = 131352.5
$ws.Range("M34").Value = -131150.5
$ws.Range("H132").Value = 20523298
$ws.Range("I132").Value = 27790248
$ws.Range("K132").Value = 83370744
$ws.Range("M132").Value = -83368214
$ws.Range("H139").Value = 94000
$ws.Range("J139").Value = 110000
$ws.Range("L139").Value = 110000
$ws.Range("N139").Value = -120280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 115679.336
$ws.Range("I109").Value = 169552.33
$ws.Range("J109").Value = 7933.3335
$ws.Range("K109").Value = 508656.99
$ws.Range("L109").Value = 23800.0005
$ws.Range("M109").Value = -507616.99
$ws.Range("N109").Value = -25880.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 592.625
$ws.Range("I31").Value = 592.625
$ws.Range("K31").Value = 592.625
$ws.Range("M31").Value = -300.625
$ws.Range("H37").Value = 592.625
$ws.Range("I37").Value = 592.625
$ws.Range("K37").Value = 592.625
$ws.Range("M37").Value = -315.625
$ws.Range("H39").Value = 36461.5
$ws.Range("J39").Value = 42998.668
$ws.Range("L39").Value = 42998.668
$ws.Range("N39").Value = -44062.668
$ws.Range("H80").Value = 252521.4
$ws.Range("I80").Value = 359567.44
$ws.Range("J80").Value = 2747.3333
$ws.Range("K80").Value = 359567.44
$ws.Range("L80").Value = 2747.3333
$ws.Range("M80").Value = -358569.44
$ws.Range("N80").Value = -4743.3333
$ws.Range("H83").Value = 252521.4
$ws.Range("I83").Value = 359567.44
$ws.Range("J83").Value = 2747.3333
$ws.Range("K83").Value = 1797837.2
$ws.Range("L83").Value = 13736.6665
$ws.Range("M83").Value = -1792845.2
$ws.Range("N83").Value = -23720.6665
$ws.Range("H95").Value = 87275.2
$ws.Range("J95").Value = 87275.2
$ws.Range("L95").Value = 87275.2
$ws.Range("N95").Value = -92767.2
$ws.Range("H96").Value = 34081.332
$ws.Range("J96").Value = 34081.332
$ws.Range("L96").Value = 34081.332
$ws.Range("N96").Value = -39573.332
$ws.Range("H126").Value = 1044737.9
$ws.Range("J126").Value = 3845
$ws.Range("L126").Value = 11535
$ws.Range("N126").Value = -16475

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1180.1428
$ws.Range("J46").Value = 1387.75
$ws.Range("L46").Value = 1387.75
$ws.Range("N46").Value = -1763.75
$ws.Range("H55").Value = 31250290
$ws.Range("I55").Value = 440
$ws.Range("J55").Value = 50000200
$ws.Range("K55").Value = 440
$ws.Range("L55").Value = 50000200
$ws.Range("M55").Value = -267
$ws.Range("N55").Value = -50000546
$ws.Range("H82").Value = 1032.1177
$ws.Range("J82").Value = 1235
$ws.Range("L82").Value = 1235
$ws.Range("N82").Value = -1957
$ws.Range("H85").Value = 1032.1177
$ws.Range("J85").Value = 1235
$ws.Range("L85").Value = 1235
$ws.Range("N85").Value = -3731
$ws.Range("H93").Value = 1152.1538
$ws.Range("I93").Value = 452.77777
$ws.Range("K93").Value = 452.77777
$ws.Range("M93").Value = 795.2222300000001
$ws.Range("H132").Value = 1453564.6
$ws.Range("I132").Value = 2322204.2
$ws.Range("J132").Value = 5832.1113
$ws.Range("K132").Value = 6966612.600000001
$ws.Range("L132").Value = 17496.3339
$ws.Range("M132").Value = -6964082.600000001
$ws.Range("N132").Value = -22556.3339
$ws.Range("H136").Value = 42057.938
$ws.Range("I136").Value = 1823.2632
$ws.Range("K136").Value = 5469.7896
$ws.Range("M136").Value = -2919.7896

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3432.1667
$ws.Range("I62").Value = 3166.6667
$ws.Range("J62").Value = 3697.6667
$ws.Range("K62").Value = 3166.6667
$ws.Range("L62").Value = 3697.6667
$ws.Range("M62").Value = -2542.6667
$ws.Range("N62").Value = -4945.6667
$ws.Range("H65").Value = 3432.1667
$ws.Range("I65").Value = 3166.6667
$ws.Range("J65").Value = 3697.6667
$ws.Range("K65").Value = 15833.3335
$ws.Range("L65").Value = 18488.3335
$ws.Range("M65").Value = -12713.3335
$ws.Range("N65").Value = -24728.3335
$ws.Range("H81").Value = 2581.8
$ws.Range("I81").Value = 1448.091
$ws.Range("K81").Value = 2896.182
$ws.Range("M81").Value = -1835.182
$ws.Range("H84").Value = 2581.8
$ws.Range("I84").Value = 1448.091
$ws.Range("K84").Value = 14480.91
$ws.Range("M84").Value = -9176.91
$ws.Range("H96").Value = 1893.75
$ws.Range("I96").Value = 1790
$ws.Range("K96").Value = 1790
$ws.Range("M96").Value = -417
$ws.Range("H100").Value = 1748.3334
$ws.Range("I100").Value = 1466.4117
$ws.Range("K100").Value = 2932.8234
$ws.Range("M100").Value = -2391.8234
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H136").Value = 1651190.9
$ws.Range("I136").Value = 2008226.5
$ws.Range("J136").Value = 44530.5
$ws.Range("K136").Value = 6024679.5
$ws.Range("L136").Value = 133591.5
$ws.Range("M136").Value = -6022129.5
$ws.Range("N136").Value = -138691.5
